# Add results for the push-up event (row 11) and the next week's partial
# entry (row 12), then move the active selection like the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: full week of results -----------------------------------------
# Copy the date-formatted style (s="1") from A2 down to A11/A12 first, then
# overwrite the values so the numeric cell format ("m/d/yyyy") is preserved
# without Excel minting a brand-new style index.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null

$ws.Range("A11").Value = 44590
$ws.Range("B11").Value = 51
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 123
$ws.Range("F11").Value = 35
$ws.Range("G11").Value = 38
$ws.Range("H11").Value = 60
$ws.Range("I11").Value = 25
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("T11").Value = 100
$ws.Range("U11").Value = 65
$ws.Range("V11").Value = 80

$ws.Range("W11").Formula = "=B11+C11+D11"
$ws.Range("X11").Formula = "=E11+F11+G11"
$ws.Range("Y11").Formula = "=H11+I11+J11"
$ws.Range("Z11").Formula = "=K11+L11+M11"
$ws.Range("AA11").Formula = "=N11+O11+P11"
$ws.Range("AB11").Formula = "=Q11+R11+S11"
$ws.Range("AC11").Formula = "=T11+U11+V11"

# --- Row 12: next week, only partial results entered so far ---------------
$ws.Range("A12").Value = 44611
$ws.Range("E12").Value = 200
$ws.Range("F12").Value = 350
$ws.Range("G12").Value = 455

# --- Leave the selection where the author ended up -------------------------
$ws.Range("H15").Select() | Out-Null
